$wb = $excel.ActiveWorkbook

# Update "Last Updated" timestamp on the Metadata sheet
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("A2").Value = "05 Nov 2025, 11:18 AM"

# Update the "1 Year" column (F) values on the Industry Analysis sheet
$ws = $wb.Worksheets.Item("Industry Analysis")

$ws.Cells.Item(2, 6).Value = 21.3
$ws.Cells.Item(3, 6).Value = -4.3927
$ws.Cells.Item(4, 6).Value = 35.9445
$ws.Cells.Item(5, 6).Value = -51.0482
$ws.Cells.Item(6, 6).Value = 57.2275
$ws.Cells.Item(7, 6).Value = -9.640700000000001
$ws.Cells.Item(8, 6).Value = -6.1449
$ws.Cells.Item(9, 6).Value = 36.9733
$ws.Cells.Item(10, 6).Value = -4.7026
$ws.Cells.Item(11, 6).Value = 46.5317
$ws.Cells.Item(12, 6).Value = -2.102
$ws.Cells.Item(13, 6).Value = 17.4681
$ws.Cells.Item(14, 6).Value = -33.0245
$ws.Cells.Item(15, 6).Value = 1.0205
$ws.Cells.Item(16, 6).Value = 2.0426
$ws.Cells.Item(17, 6).Value = -16.2411
$ws.Cells.Item(18, 6).Value = 7.4627
$ws.Cells.Item(19, 6).Value = -25.798
$ws.Cells.Item(20, 6).Value = 47.7485
$ws.Cells.Item(21, 6).Value = 19.5587
$ws.Cells.Item(22, 6).Value = 76.5603
$ws.Cells.Item(23, 6).Value = -54.2675
$ws.Cells.Item(24, 6).Value = -0.8811
$ws.Cells.Item(25, 6).Value = 4.8518
$ws.Cells.Item(26, 6).Value = 3.6831
$ws.Cells.Item(27, 6).Value = -34.0874
$ws.Cells.Item(28, 6).Value = -11.9893
$ws.Cells.Item(29, 6).Value = -12.994
$ws.Cells.Item(30, 6).Value = 25.5415
$ws.Cells.Item(31, 6).Value = 56.5088
$ws.Cells.Item(32, 6).Value = 2.0908
$ws.Cells.Item(33, 6).Value = -4.7193
$ws.Cells.Item(34, 6).Value = 22.8807
$ws.Cells.Item(35, 6).Value = 5.3359
$ws.Cells.Item(36, 6).Value = -5.1995
$ws.Cells.Item(37, 6).Value = -5.6238
$ws.Cells.Item(38, 6).Value = -22.595
$ws.Cells.Item(39, 6).Value = 10.8405
$ws.Cells.Item(40, 6).Value = -7.5963
$ws.Cells.Item(41, 6).Value = -4.552
$ws.Cells.Item(42, 6).Value = 22.3098
$ws.Cells.Item(43, 6).Value = 14.0694
$ws.Cells.Item(44, 6).Value = -9.6066
$ws.Cells.Item(45, 6).Value = 27.639
$ws.Cells.Item(46, 6).Value = -6.3484
$ws.Cells.Item(47, 6).Value = -40.5302
$ws.Cells.Item(48, 6).Value = -29.7988
$ws.Cells.Item(49, 6).Value = -24.0791
$ws.Cells.Item(50, 6).Value = -49.1803
$ws.Cells.Item(51, 6).Value = -51.6023
$ws.Cells.Item(52, 6).Value = -34.4756
$ws.Cells.Item(53, 6).Value = -11.5478
$ws.Cells.Item(54, 6).Value = -2.3796
$ws.Cells.Item(55, 6).Value = -15.4382
$ws.Cells.Item(56, 6).Value = -27.6987
$ws.Cells.Item(57, 6).Value = -27.1559
$ws.Cells.Item(58, 6).Value = -2.1585
$ws.Cells.Item(59, 6).Value = -23.0964
$ws.Cells.Item(60, 6).Value = -13.3217
$ws.Cells.Item(61, 6).Value = -8.1496
$ws.Cells.Item(62, 6).Value = -16.0695
$ws.Cells.Item(63, 6).Value = -12.5465
$ws.Cells.Item(64, 6).Value = 47.7264
$ws.Cells.Item(65, 6).Value = -42.4232
$ws.Cells.Item(66, 6).Value = 11.3291
$ws.Cells.Item(67, 6).Value = 14.3746
$ws.Cells.Item(68, 6).Value = 32.6702
$ws.Cells.Item(69, 6).Value = -17.0097
$ws.Cells.Item(70, 6).Value = -13.5162
$ws.Cells.Item(71, 6).Value = 11.4259
$ws.Cells.Item(72, 6).Value = 2.6754
$ws.Cells.Item(73, 6).Value = -11.1574
$ws.Cells.Item(74, 6).Value = -13.2502
$ws.Cells.Item(75, 6).Value = 24.7078
$ws.Cells.Item(76, 6).Value = 53.3554